$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3095.037
$ws.Range("I32").Value = 2033.1666
$ws.Range("J32").Value = 3398.4285
$ws.Range("K32").Value = 2033.1666
$ws.Range("L32").Value = 3398.4285
$ws.Range("M32").Value = -1707.1666
$ws.Range("N32").Value = -4050.4285

$ws.Range("H43").Value = 1224210.5
$ws.Range("J43").Value = 1224210.5
$ws.Range("L43").Value = 1224210.5
$ws.Range("N43").Value = -1224348.5

$ws.Range("H57").Value = 199999
$ws.Range("J57").Value = 199999
$ws.Range("L57").Value = 599997
$ws.Range("N57").Value = -600995

$ws.Range("H62").Value = 1573.75
$ws.Range("I62").Value = 1573.75
$ws.Range("K62").Value = 1573.75
$ws.Range("M62").Value = -949.75

$ws.Range("H65").Value = 1573.75
$ws.Range("I65").Value = 1573.75
$ws.Range("K65").Value = 7868.75
$ws.Range("M65").Value = -4748.75

$ws.Range("H113").Value = 6930.4346
$ws.Range("I113").Value = 4431.25
$ws.Range("J113").Value = 12642.857
$ws.Range("K113").Value = 4431.25
$ws.Range("L113").Value = 12642.857
$ws.Range("M113").Value = -1177.25
$ws.Range("N113").Value = -19150.857

$ws.Range("H138").Value = 100001640
$ws.Range("I138").Value = 1820.6666
$ws.Range("J138").Value = 1000000000
$ws.Range("K138").Value = 5461.9998
$ws.Range("L138").Value = 3000000000
$ws.Range("M138").Value = -321.9997999999996
$ws.Range("N138").Value = -3000010280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3200.4941
$ws.Range("I32").Value = 3258.1643
$ws.Range("K32").Value = 3258.1643
$ws.Range("M32").Value = -2971.1643

$ws.Range("H43").Value = 250036740
$ws.Range("J43").Value = 48990
$ws.Range("L43").Value = 48990
$ws.Range("N43").Value = -49616

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 32117.445
$ws.Range("J41").Value = 39142.57
$ws.Range("L41").Value = 39142.57
$ws.Range("N41").Value = -39998.57

$ws.Range("H50").Value = 33099.8
$ws.Range("J50").Value = 33099.8
$ws.Range("L50").Value = 33099.8
$ws.Range("N50").Value = -34349.8

$ws.Range("H51").Value = 25000
$ws.Range("I51").Value = 25000
$ws.Range("K51").Value = 25000
$ws.Range("M51").Value = -24264

$ws.Range("H60").Value = 3887.4
$ws.Range("I60").Value = 3887.4
$ws.Range("K60").Value = 3887.4
$ws.Range("M60").Value = -3376.4

$ws.Range("H61").Value = 25000
$ws.Range("I61").Value = 25000
$ws.Range("K61").Value = 25000
$ws.Range("M61").Value = -24652

$ws.Range("H105").Value = 8672.454
$ws.Range("I105").Value = 9249.700000000001
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 9249.700000000001
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -7502.700000000001
$ws.Range("N105").Value = -6394

$ws.Range("H131").Value = 33370
$ws.Range("J131").Value = 37644
$ws.Range("L131").Value = 37644
$ws.Range("N131").Value = -47724

$ws.Range("H132").Value = 1931.0256
$ws.Range("I132").Value = 2043.5938
$ws.Range("K132").Value = 6130.7814
$ws.Range("M132").Value = -3600.7814

$ws.Range("H138").Value = 199996.67
$ws.Range("J138").Value = 199996.67
$ws.Range("L138").Value = 199996.67
$ws.Range("N138").Value = -210276.67

$ws.Range("H141").Value = 99120.27
$ws.Range("J141").Value = 123165.375
$ws.Range("L141").Value = 123165.375
$ws.Range("N141").Value = -133525.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1806.1072
$ws.Range("I5").Value = 1060.9375
$ws.Range("K5").Value = 3182.8125
$ws.Range("M5").Value = -3070.8125

$ws.Range("H31").Value = 3637.2727
$ws.Range("I31").Value = 3637.2727
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 10911.8181
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -10623.8181
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 102153.695
$ws.Range("J37").Value = 102153.695
$ws.Range("L37").Value = 306461.085
$ws.Range("N37").Value = -306685.085

$ws.Range("H107").Value = 366.05
$ws.Range("J107").Value = 357.06668
$ws.Range("L107").Value = 1071.20004
$ws.Range("N107").Value = -4911.20004

$ws.Range("H132").Value = 3522.4119
$ws.Range("I132").Value = 2133.5715
$ws.Range("J132").Value = 4494.6
$ws.Range("K132").Value = 19202.1435
$ws.Range("L132").Value = 40451.4
$ws.Range("M132").Value = -16672.1435
$ws.Range("N132").Value = -45511.4

$ws.Range("H135").Value = 1806.1072
$ws.Range("I135").Value = 1060.9375
$ws.Range("K135").Value = 9548.4375
$ws.Range("M135").Value = -7013.4375

$ws.Range("H137").Value = 1699.6
$ws.Range("I137").Value = 999.4286
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 2998.2858
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = 2101.7142
$ws.Range("N137").Value = -20199.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 27371
$ws.Range("I93").Value = 19500
$ws.Range("J93").Value = 29994.666
$ws.Range("K93").Value = 19500
$ws.Range("L93").Value = 29994.666
$ws.Range("M93").Value = -17628
$ws.Range("N93").Value = -33738.666

$ws.Range("H129").Value = 29102.25
$ws.Range("J129").Value = 39994.5
$ws.Range("L129").Value = 39994.5
$ws.Range("N129").Value = -49994.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1928.9
$ws.Range("J16").Value = 1510.25
$ws.Range("L16").Value = 1510.25
$ws.Range("N16").Value = -1850.25

$ws.Range("H22").Value = 1650.4
$ws.Range("I22").Value = 1278.2222
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 1278.2222
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -983.2221999999999
$ws.Range("N22").Value = -5590

$ws.Range("H27").Value = 1650.4
$ws.Range("I27").Value = 1278.2222
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 1278.2222
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1171.2222
$ws.Range("N27").Value = -5214

$ws.Range("H122").Value = 3263.4
$ws.Range("I122").Value = 2905.5
$ws.Range("J122").Value = 4695
$ws.Range("K122").Value = 8716.5
$ws.Range("L122").Value = 14085
$ws.Range("M122").Value = -6266.5
$ws.Range("N122").Value = -18985

$ws.Range("H131").Value = 65998.5
$ws.Range("J131").Value = 75664.664
$ws.Range("L131").Value = 75664.664
$ws.Range("N131").Value = -85744.664

$ws.Range("H132").Value = 10870.5
$ws.Range("I132").Value = 3783.7896
$ws.Range("K132").Value = 11351.3688
$ws.Range("M132").Value = -8821.3688

$ws.Range("H136").Value = 2504.0605
$ws.Range("I136").Value = 2313.7083
$ws.Range("J136").Value = 3011.6667
$ws.Range("K136").Value = 6941.124899999999
$ws.Range("L136").Value = 9035.000100000001
$ws.Range("M136").Value = -4391.124899999999
$ws.Range("N136").Value = -14135.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 149499.5
$ws.Range("J131").Value = 149499.5
$ws.Range("L131").Value = 149499.5
$ws.Range("N131").Value = -159579.5

$ws.Range("H132").Value = 2449.6775
$ws.Range("I132").Value = 2227.45
$ws.Range("K132").Value = 6682.349999999999
$ws.Range("M132").Value = -4152.349999999999

$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
